# T460 - update the June 2022 "Foaie de parcurs" travel log for B 151 VGT
# (Alex Bora): corrected daily km / destination / purpose entries and the
# resulting running totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Starting odometer reading
$ws.Range("B12").Value = 307760

# Day 2 (row 15)
$ws.Range("B15").Value = 85
$ws.Range("C15").Value = "Cluj-Apahida"
$ws.Range("D15").Value = "Interes Serviciu"

# Day 3 (row 16)
$ws.Range("B16").Value = 47
$ws.Range("C16").Value = "Cluj-Cluj"
$ws.Range("D16").Value = "Interes Serviciu"

# Day 6 (row 19)
$ws.Range("B19").Value = 92
$ws.Range("C19").Value = "Cluj-Bontida"
$ws.Range("D19").Value = "Interes Serviciu"

# Day 7 (row 20)
$ws.Range("B20").Value = 47
$ws.Range("C20").Value = "Cluj-Cluj"
$ws.Range("D20").Value = "Interes Serviciu"

# Day 8 (row 21)
$ws.Range("B21").Value = 257
$ws.Range("C21").Value = "Cluj-Bistrita"
$ws.Range("D21").Value = "Interes Serviciu"

# Day 9 (row 22)
$ws.Range("B22").Value = 30
$ws.Range("C22").Value = "Acasa-Birou"
$ws.Range("D22").Value = " "

# Day 10 (row 23)
$ws.Range("B23").Value = 30
$ws.Range("C23").Value = "Acasa-Birou"
$ws.Range("D23").Value = " "

# Day 14 (row 27)
$ws.Range("B27").Value = 421
$ws.Range("C27").Value = "Cluj-Satu-Mare"
$ws.Range("D27").Value = "Interes Serviciu"

# Day 15 (row 28)
$ws.Range("B28").Value = 152
$ws.Range("C28").Value = "Cluj-Cmp. Turzii"
$ws.Range("D28").Value = "Interes Serviciu"

# Day 16 (row 29)
$ws.Range("B29").Value = 30
$ws.Range("C29").Value = "Acasa-Birou"
$ws.Range("D29").Value = " "

# Day 17 (row 30)
$ws.Range("B30").Value = 121
$ws.Range("C30").Value = "Cluj-Turda"
$ws.Range("D30").Value = "Interes Serviciu"

# Day 20 (row 33)
$ws.Range("B33").Value = 30
$ws.Range("C33").Value = "Acasa-Birou"
$ws.Range("D33").Value = " "

# Day 21 (row 34)
$ws.Range("B34").Value = 92
$ws.Range("C34").Value = "Cluj-Bontida"
$ws.Range("D34").Value = "Interes Serviciu"

# Day 23 (row 36)
$ws.Range("B36").Value = 156
$ws.Range("C36").Value = "Cluj-Zalau"
$ws.Range("D36").Value = "Interes Serviciu"

# Day 24 (row 37)
$ws.Range("B37").Value = 47
$ws.Range("C37").Value = "Cluj-Cluj"
$ws.Range("D37").Value = "Interes Serviciu"

# Day 27 (row 40)
$ws.Range("B40").Value = 121
$ws.Range("C40").Value = "Cluj-Turda"
$ws.Range("D40").Value = "Interes Serviciu"

# Day 28 (row 41)
$ws.Range("B41").Value = 85
$ws.Range("C41").Value = "Cluj-Apahida"
$ws.Range("D41").Value = "Interes Serviciu"

# Day 29 (row 42)
$ws.Range("B42").Value = 30
$ws.Range("C42").Value = "Acasa-Birou"
$ws.Range("D42").Value = " "

# Day 30 (row 43)
$ws.Range("B43").Value = 30
$ws.Range("C43").Value = "Acasa-Birou"
$ws.Range("D43").Value = " "

# Monthly totals
$ws.Range("B44").Value = 2259
$ws.Range("B45").Value = 310019
